$d = $word.ActiveDocument

# --- Simple text replacements (no tab runs involved) ---

$found0 = $d.Content.Find.Execute("October 28, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "November 3, 2024", 2)
if (-not $found0) { Write-Output "REPLACEMENT 0 FAILED" }

$found1 = $d.Content.Find.Execute("Capacitive touch screens may be one of the most amazing pieces of technology ever, just some magical layers of metal, plastic, and glass that all fit together to make something that can sense a finger. Capacitors are probably even more important as they keep your TV on during a brief fraction of a second power outage, or make a clock signal inside your computer.  Things you can learn about these amazing technologies are the history of capacitors, how they function, the pros and cons of capacitive touch screens, and how they work.", $true, $false, $false, $false, $false, $true, 1, $false, "Capacitive touch screens may be one of the most important pieces of technology ever— just layers of metal, plastic, and glass that all fit together just right to make something that can detect a finger. Not only are capacitive touch screens amazing inventions, but capacitors are also extremely important because they do other things like create the clock signal for all your electronic devices.  Some things you can learn about these revolutionary technologies are the history of capacitors, the pros and cons of capacitive touch screens, and how they both work.", 2)
if (-not $found1) { Write-Output "REPLACEMENT 1 FAILED" }

$found2 = $d.Content.Find.Execute("Some amazing things you can learn about capacitors is their history. The first person to ever build a capacitor was Pieter van Musschenbroek when he invented the Leyden Jar (Brain & Pollette, 2021). The Leyden Jar was a glass jar filled halfway with water, lined inside and out with metal foil, and a metal wire going through the top connected to a power supply. Later, Benjamin Franklin invented a newer variation that was flat, called the Franklin Square. The first unit ever used to measure the capacitance of a capacitor was the Farad invented by Micheal Faraday, which is still used today. A one Farad capacitor can store 1 coulomb of charge at one volt (that’s 6.25 billion electrons.)", $true, $false, $false, $false, $false, $true, 1, $false, "One key thing to know about capacitors is how they function. Basically, a capacitor is a device that is used to store electrical charge and energy when used with a DC current. Capacitors are made up of two plates separated by either a dielectric (a material that does not conduct electricity) or a vacuum (Ling et al., 2016). When connected to a DC current with one plate positively charged and the other negatively charged, almost no electrons can make it across the gap between the two plates. The negatively and positively charged plates create an electric field between them, which can be detected in some cases (Ling et al., 2016).", 2)
if (-not $found2) { Write-Output "REPLACEMENT 2 FAILED" }

$found3 = $d.Content.Find.Execute("The most important thing to know about capacitors is how they work. Basically, a capacitor is a device that is used to store electrical charge and energy. Capacitors are made up of two plates separated by either a dielectric or a vacuum (Ling et al., 2016). When connected to a DC current with one plate positively charged and the other negatively charged, almost no electrons can make it across the gap between the two plates. The negatively and positively charged plates create an electric field between them, which can be detected in some cases.", $true, $false, $false, $false, $false, $true, 1, $false, "Now that we know about how capacitors function, it is time to learn about their history. The history of the capacitor starts out with Pieter van Musschenbroek, who invented the first capacitor called the Leyden Jar (Brain & Pollette, 2021). The Leyden Jar was a glass jar filled halfway with water that was lined inside and out with metal foil and had a metal wire going through the top that was connected to a power supply. Even though the Leyden Jar was not useful at the time of its invention, it was still an important step toward the capacitors we use today. Later, Benjamin Franklin invented a newer variation that was flat, called the Franklin Square. The first unit ever used to measure the capacitance (the amount of energy a capacitor or circuit can store) of a capacitor was the farad invented by Micheal Faraday, which is still used today. A one farad capacitor can store one coulomb (6.25 billion billion electrons) of charge at one volt (Brain & Pollette, 2021).", 2)
if (-not $found3) { Write-Output "REPLACEMENT 3 FAILED" }

# --- Paragraph-level replacements for paragraphs that start with a <w:tab/> run, ---
# --- using InsertXML on the whole paragraph Range so the <w:tab/> element survives ---
# --- (plain Find/Replace flattens the leading tab into literal text). ---

$p4 = $d.Paragraphs.Item(13)
$rng4 = $p4.Range
if ($rng4.Text.IndexOf("One important aspect to learn about capacitive touch screens are their pros, cons, and capabilities. Capacitive touch screens have some key advantages over other types of touch screens. These advantages include smoother and faster scrolling, the ability to use any material as the top layer, and full support for multitouch (Nam et al., 2021). Because capacitive touch screens can be made of almost any material, they can be made to be more durable to chemicals and forces by using glass or plastic (Barrett & Omote, 2010). Even though capacitive touch screens have their advantages, they also have their downsides. Some of these downsides include water being able to trigger touches, not being able to use gloves, not being able to work with objects not designed to be used with a capacitive touch screen, and being limited to a smaller size (Glinpu, 2023; Nam et al., 2021).") -lt 0) { Write-Output "PARAGRAPH 4 TEXT MISMATCH" }
$rng4.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"normal1`"/><w:spacing w:lineRule=`"auto`" w:line=`"480`"/><w:ind w:hanging=`"0`" w:left=`"0`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>Capacitive touch screens have many pros, cons, and capabilities. Some of the capabilities of capacitive touch screens include smoother and faster scrolling, the ability to use any material as the top layer, and full support for multitouch (Nam et al., 2021). Because capacitive touch screens can be made of almost any material, they can be made to be more durable against chemicals and forces by using materials like glass or plastic (Barrett &amp; Omote, 2010). Even though capacitive touch screens have their advantages, they also have their downsides. Some of these downsides include water being able to trigger touches, not being able to use gloves, being limited to objects designed specifically for this kind of touch screen, and being limited to a smaller size (Glinpu, 2023; Nam et al., 2021).</w:t></w:r></w:p>")

$p5 = $d.Paragraphs.Item(14)
$rng5 = $p5.Range
if ($rng5.Text.IndexOf("Most people use their fingers to control touch screens, but how does this work? A capacitive touch screen is made up of multiple layers. There are two key layers when it comes to detecting a finger, the top layer usually made up of glass or plastic, and the layer below that is made up of indium tin oxide or ITO (Glinpu, 2023). ITO is one plate of a capacitor, your finger is the other, and the glass or plastic is the dielectric in between. The ITO layer is powered by an AC current which unlike a DC current can easily flow through a capacitor but still makes an electric field. Because your body is always connected to an electrical ground, the voltage and wattage are safe enough to flow through your body with no harm (Saini, 2011).") -lt 0) { Write-Output "PARAGRAPH 5 TEXT MISMATCH" }
$rng5.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"normal1`"/><w:spacing w:lineRule=`"auto`" w:line=`"480`"/><w:ind w:hanging=`"0`" w:left=`"0`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>Many people use their fingers to control capacitive touch screens, but how do they work? A capacitive touch screen is made up of multiple layers. There are two key layers when it comes to detecting a finger—the top layer usually made up of glass or plastic, and the layer below that is made up of indium tin oxide (ITO) (Glinpu, 2023). ITO is one plate of a capacitor, your finger is the other, and the glass or plastic is the dielectric in between. The ITO layer is powered by an AC current that, unlike a DC current, can easily flow through a capacitor but still makes an electric field. Because your body is always connected to an electrical ground, the voltage and wattage are safe enough to flow through your body with no harm (Saini, 2011).</w:t></w:r></w:p>")

$p6 = $d.Paragraphs.Item(15)
$rng6 = $p6.Range
if ($rng6.Text.IndexOf("The most important thing to understand about capacitive touch screens is how they detect fingers. There are two main ways capacitive touch screens detect a finger, surface capacitive touch screens and projected capacitive touch screens. Surface capacitive touch screens are made up of one layer of conductive material with all 4 corners connected to a synchronized AC current (Nam et al., 2021). When your finger comes in contact with the touch screen, it creates a current difference. Because the four corners are different distances from the point of contact, the controller can determine the location of the finger based on the current difference in each corner.") -lt 0) { Write-Output "PARAGRAPH 6 TEXT MISMATCH" }
$rng6.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"normal1`"/><w:spacing w:lineRule=`"auto`" w:line=`"480`"/><w:ind w:hanging=`"0`" w:left=`"0`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>There are multiple different types of capacitive touch screens, each with their own purpose. There are two main ways capacitive touch screens detect a finger—surface-capacitive touch screens, and projected-capacitive touch screens. Surface-capacitive touch screens are made up of one layer of conductive material with all four corners connected to a synchronized AC current (Nam et al., 2021). When your finger comes in contact with the touch screen, it creates a current difference. Because the four corners are different distances from the point of contact, the controller can determine the location of the finger based on the current difference in each corner (Nam et al., 2021).</w:t></w:r></w:p>")

$p7 = $d.Paragraphs.Item(16)
$rng7 = $p7.Range
if ($rng7.Text.IndexOf("The other kind of capacitive touch screen is the projected capacitive touch screen. Projected capacitive touch screens consist of two layers of lines making an intersecting pattern, usually a grid of rows and columns (Nam et al., 2021). Projected capacitive touch screens have two subtypes, self capacitive touch screens and mutual capacitive touch screens. Self capacitive touch screens work by scanning all the rows and columns, determining coordinates from the columns and rows with a current difference (Barrett & Omote, 2010). This causes self capacitive touch screens to have ghost points where you have two points one position off in each axis, making it impossible to tell where the points are. Mutual capacitive touch screens detect at each intersection instead of at each row or column (Barrett & Omote, 2010. This gives them full multitouch without ghost points.") -lt 0) { Write-Output "PARAGRAPH 7 TEXT MISMATCH" }
$rng7.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"normal1`"/><w:spacing w:lineRule=`"auto`" w:line=`"480`"/><w:ind w:hanging=`"0`" w:left=`"0`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>The other kind of capacitive touch screen is the projected-capacitive touch screen. Projected-capacitive touch screens consist of two layers of lines making an intersecting pattern, usually a grid of rows and columns (Nam et al., 2021). Projected-capacitive touch screens have two subtypes—the self-capacitance and mutual-capacitance architectures. The self-capacitance architecture works by scanning all the rows and columns, then determining coordinates from the columns and rows with a current difference (Barrett &amp; Omote, 2010). This causes the self-capacitance architecture to have ghost points when you have two points one position off in each axis, making it impossible to tell where the points are. The mutual-capacitance architecture detects at each intersection instead of at each row or column; this gives them full multitouch without ghost points (Barrett &amp; Omote, 2010).</w:t></w:r></w:p>")

$p8 = $d.Paragraphs.Item(17)
$rng8 = $p8.Range
if ($rng8.Text.IndexOf("Capacitive touch screens have improved technology in the modern era by allowing multitouch support and scrolling smoother and faster. There are many different types of capacitive touch screens. My experiment will focus on capacitive buttons, which are similar to their touch screen counterpart in the sense that they both detect a finger in the same way. By testing different patterns of capacitive buttons, I want to test which pattern works through the most layers of plastic or other dielectric. After doing research related to this topic, I hypothesize that all the patterns will work about the same because the button will always make a capacitor with your finger regardless of the pattern, which is all that is necessary to detect a finger.") -lt 0) { Write-Output "PARAGRAPH 8 TEXT MISMATCH" }
$rng8.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:pPr><w:pStyle w:val=`"normal1`"/><w:spacing w:lineRule=`"auto`" w:line=`"480`"/><w:ind w:hanging=`"0`" w:left=`"0`"/><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`"/><w:b/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:sz w:val=`"24`"/><w:szCs w:val=`"24`"/></w:rPr><w:tab/><w:t>Capacitive touch screens have improved technology in the modern era by allowing multitouch support and by making scrolling smoother and faster. There are many different types of capacitive touch screens. My experiment will focus on capacitive buttons, which are similar to their touch screen counterpart in the sense that they both detect a finger using the same method. By testing different patterns of capacitive buttons, I want to test which pattern works through the most layers of plastic. After doing research related to this topic, I hypothesize that all the patterns will work about the same because the metal layer of each button will always form a capacitor with your finger regardless of the pattern, which is all that is necessary to detect a finger.</w:t></w:r><w:r><w:br w:type=`"page`"/></w:r></w:p>")

Write-Output "Done"